$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
Write-Output ($ws.Range("A1").Value)
Write-Output ($ws.Range("A1").Value2)
Write-Output ($ws.Range("A1").Text)
